$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = "new changes done from f2"
$ws.Range("D8").Select()
